$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") computed from regenerated save_data (K = walks instead of old Strike# count)
$values = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 4
    6  = 1
    7  = 4
    8  = 2
    9  = 2
    10 = 7
    11 = 3
    12 = 6
    13 = 8
    14 = 4
    15 = 2
    16 = 4
    17 = 3
    18 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
